$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A (Name) for the new "Dee" investigation rows first
$ws.Range("A2").Value = "Dee"
$ws.Range("A3").Value = "Dee"
$ws.Range("A4").Value = "Dee"

# Fill column C (Avatar) next
$ws.Range("C2").Value = "Dee-Thinking2"
$ws.Range("C3").Value = "Dee-Determined"
$ws.Range("C4").Value = "Dee-Thinking2"

# Fill the BGM column for row 2
$ws.Range("L2").Value = "Dee-Thinking"

# Fill column B (Dialogue) along with the new BGImage tag (F2)
$ws.Range("B2").Value = "On the outer side of the door panel, there are bloodstains resembling handprints."
$ws.Range("B3").Value = "The prints are smudged with elongated streaks, suggesting they were dragged up and down."
$ws.Range("F2").Value = "Suspicious"
$ws.Range("B4").Value = "There are also traces of blood on the door knocker of the manor" + [char]8217 + "s main gate" + [char]8212 + [char]8212 + "completely dried and coagulated."

# Row 5 - remove old D5/E5 values (no longer carried for this row)
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()

# Row 6 - remove old D6/E6 values
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()

# Update row heights to match wrapped text content
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 51

# Update selected cell (UI state) to match the new selection
$ws.Range("B10").Select()
